# Updates odds/score values on Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("I2").Value = 2.35
$ws.Range("L2").Value = 3.25
$ws.Range("AE2").Value = 15

# Row 3
$ws.Range("G3").Value = 2.9
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 2.2
$ws.Range("L3").Value = 2.88
$ws.Range("Y3").Value = 1.33
$ws.Range("Z3").Value = 3.25
$ws.Range("AA3").Value = 1.57
$ws.Range("AB3").Value = 2.25
$ws.Range("AC3").Value = 12
$ws.Range("AD3").Value = 17
$ws.Range("AH3").Value = 26
$ws.Range("AI3").Value = 13
$ws.Range("AJ3").Value = 7
$ws.Range("AM3").Value = 10
$ws.Range("AN3").Value = 12
$ws.Range("AP3").Value = 21
$ws.Range("AS3").Value = 126

# Row 4
$ws.Range("U4").Value = 2.03
$ws.Range("V4").Value = 1.83

# Row 8
$ws.Range("H8").Value = 3.4
$ws.Range("I8").Value = 3.8
$ws.Range("K8").Value = 2.05
$ws.Range("M8").Value = 1.06
$ws.Range("N8").Value = 8
$ws.Range("S8").Value = 2.05
$ws.Range("T8").Value = 1.75
$ws.Range("AA8").Value = 1.91
$ws.Range("AB8").Value = 1.8
$ws.Range("AC8").Value = 6.5
$ws.Range("AI8").Value = 9
$ws.Range("AO8").Value = 15

$wb.Save()
